$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value = "Move Robot2 to location (2, 8) and remove the toolkit."
$ws.Range("B1").Value = "['Robot2']"
$ws.Range("E1").Value = "(2, 8)"

# Row 2
$ws.Range("A2").Value = "Move Robot26 to location (11, 4) and remove the liquid spill."
$ws.Range("B2").Value = "['Robot26']"
$ws.Range("E2").Value = "(11, 4)"

# Row 3
$ws.Range("A3").Value = "Move Robot42 to location (9, 5) and remove the large debris."
$ws.Range("B3").Value = "['Robot42']"
$ws.Range("C3").Value = "['gripper']"
$ws.Range("E3").Value = "(9, 5)"

# Row 4
$ws.Range("A4").Value = "Move Robot48 to location (5, 6) and remove the dust."
$ws.Range("E4").Value = "(5, 6)"

# Row 5
$ws.Range("A5").Value = "Move Robot31 to location (9, 4) and remove the grass."
$ws.Range("B5").Value = "['Robot31']"
$ws.Range("E5").Value = "(9, 4)"

# Row 6
$ws.Range("A6").Value = "Move Robot8 to location (8, 12) and remove the small debris."
$ws.Range("B6").Value = "['Robot8']"
$ws.Range("C6").Value = "['broom']"
$ws.Range("E6").Value = "(8, 12)"

# Row 7
$ws.Range("A7").Value = "Move Robot23 to location (11, 1) and remove the vehicle."
$ws.Range("B7").Value = "['Robot23']"
$ws.Range("E7").Value = "(11, 1)"

# Row 8
$ws.Range("A8").Value = "Move Robot23 to location (12, 10) and remove the construction materials."
$ws.Range("E8").Value = "(12, 10)"

# Row 9
$ws.Range("A9").Value = "Move Robot14 to location (7, 11) and remove the tree branches."
$ws.Range("B9").Value = "['Robot14']"
$ws.Range("E9").Value = "(7, 11)"

# Row 10
$ws.Range("A10").Value = "Move Robot15 to location (5, 3) and remove the screws."
$ws.Range("E10").Value = "(5, 3)"
